$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log")

# Fill in row 8 with the new log entry
$ws.Range("B8").Value = 9417
$ws.Range("C8").Value = "31/03/2020"
$ws.Range("D8").Value = "8:45pm "
$ws.Range("E8").Value = "10:25PM "
$ws.Range("G8").Value = "Worked on the Arithmatic unit "

# Move the active cell selection to E8 (matches author's edit)
$ws.Range("E8").Select()
